# Regenerate column G ("K" = strikeouts) for each outing row, replacing the
# previous placeholder "Strike#" derived values with the recalculated K
# values, as produced by regenerating save_data (std/mean recompute + s_vals
# write upstream; only the K column changes in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (1-based, matching worksheet rows) -> new K value.
$kValues = [ordered]@{
    2  = 1
    3  = 1
    4  = 5
    5  = 2
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 2
    21 = 2
    22 = 4
    23 = 2
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 3
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 3
    47 = 1
    48 = 2
    49 = 1
    50 = 1
    51 = 2
    52 = 1
    53 = 0
    54 = 1
    55 = 2
    56 = 3
    57 = 2
    58 = 0
    59 = 1
    60 = 0
    61 = 3
    62 = 0
    63 = 2
    64 = 0
    65 = 1
    66 = 1
    67 = 1
    68 = 4
    69 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
